# Edit script: add JWT-related nuget/config rows to the "Comandos Consola" sheet.
#
# Summary of the change (see commit message: "se configuro el JWT en .net core
# API, donde se crearon nuevas clases e interfaces."):
#   - Row 7 (Scaffold-DBContext row) gains two extra notes in D7/E7 about
#     adding extra tables to the scaffolded model.
#   - Two brand new rows are inserted right after it with the nuget packages
#     needed for JWT auth (Microsoft.AspNetCore.Authentication.JwtBearer and
#     System.IdentityModel.Tokens.Jwt), highlighted with a new fill color.
#   - B6/C6/C7 (and the two new cells) get that same new highlight fill.
#   - The active selection moves down accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room: insert two new rows after the existing row 7 ------------
# (old row 8 "PROYECTO/DIAGRAMA..." header block shifts from row 8 -> row 10)
$ws.Rows.Item(8).Resize(2).Insert()

# --- 2. New text in the (still existing) row 7 ------------------------------
# (values are written E-before-D so the shared-string table ends up ordering
# the two new strings the same way the original authoring session did)
$ws.Range("E7").Value2 = 'Agregar tablas al modelo: "table1", "table2", "table3"'
$ws.Range("D7").Value2 = 'Scaffold-DBContext "Server=OFITE-GRUDE8\SQLEXPRESS;Database=VentaReal;Trusted_Connection=True;" Microsoft.EntityFrameworkCore.SqlServer -OutputDir Models -Tables "Usuario" -Force'

# --- 3. New rows 8 and 9 with the JWT nuget package references -------------
$ws.Range("C8").Value2 = 'Microsoft.AspNetCore.Authentication.JwtBearer 3.0.3'
$ws.Range("C9").Value2 = 'System.IdentityModel.Tokens.Jwt 6.6.0'

# --- 4. Apply the new highlight fill color to B6, C6, C7, C8, C9 -----------
# (solid fill equivalent to theme Accent2 tinted 80% lighter -> #FBE5D6)
$newFillColor = 14083579   # OLE BGR value for RGB FBE5D6

$ws.Range("B6").Interior.Color = $newFillColor
$ws.Range("C6").Interior.Color = $newFillColor
$ws.Range("C7").Interior.Color = $newFillColor
$ws.Range("C8").Interior.Color = $newFillColor
$ws.Range("C9").Interior.Color = $newFillColor

# --- 5. Update the active selection, matching the new layout ---------------
$ws.Range("E18").Select() | Out-Null
